# Refresh the live market-data columns (H:N) across the eight crafting-job
# "Leve profit" sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR) with the
# latest scheduled-runner pull. Plain numeric overwrites -- no formulas or
# formatting involved, only values change (a few rows also gain/lose their
# HQ-profit N cell when HQ data becomes unavailable/available).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 40
$ws.Cells.Item(40, 8).Value = 1454.8125
$ws.Cells.Item(40, 9).Value = 1616.091
$ws.Cells.Item(40, 10).Value = 1100
$ws.Cells.Item(40, 11).Value = 1616.091
$ws.Cells.Item(40, 12).Value = 1100
$ws.Cells.Item(40, 13).Value = -1441.091
$ws.Cells.Item(40, 14).Value = -1450

# Row 111
$ws.Cells.Item(111, 8).Value = 142863310
$ws.Cells.Item(111, 9).Value = 7716
$ws.Cells.Item(111, 10).Value = 250005000
$ws.Cells.Item(111, 11).Value = 23148
$ws.Cells.Item(111, 12).Value = 750015000
$ws.Cells.Item(111, 13).Value = -20081
$ws.Cells.Item(111, 14).Value = -750021134

# Row 113
$ws.Cells.Item(113, 8).Value = 6666.6294
$ws.Cells.Item(113, 9).Value = 1984.2307
$ws.Cells.Item(113, 10).Value = 11014.571
$ws.Cells.Item(113, 11).Value = 1984.2307
$ws.Cells.Item(113, 12).Value = 11014.571
$ws.Cells.Item(113, 13).Value = 1269.7693
$ws.Cells.Item(113, 14).Value = -17522.571

# Row 132
$ws.Cells.Item(132, 8).Value = 2222.2144
$ws.Cells.Item(132, 9).Value = 2423.24
$ws.Cells.Item(132, 11).Value = 7269.719999999999
$ws.Cells.Item(132, 13).Value = -4739.719999999999

$ws = $wb.Worksheets.Item("ARM")
# Row 45
$ws.Cells.Item(45, 8).Value = 2040.7727
$ws.Cells.Item(45, 9).Value = 1175.4706
$ws.Cells.Item(45, 11).Value = 1175.4706
$ws.Cells.Item(45, 13).Value = -798.4706000000001

# Row 97
$ws.Cells.Item(97, 8).Value = 1326.9375
$ws.Cells.Item(97, 9).Value = 910
$ws.Cells.Item(97, 10).Value = 2577.75
$ws.Cells.Item(97, 11).Value = 910
$ws.Cells.Item(97, 12).Value = 2577.75
$ws.Cells.Item(97, 13).Value = -414
$ws.Cells.Item(97, 14).Value = -3569.75

# Row 132
$ws.Cells.Item(132, 8).Value = 2584.1794
$ws.Cells.Item(132, 9).Value = 1286.32
$ws.Cells.Item(132, 10).Value = 4901.7856
$ws.Cells.Item(132, 11).Value = 3858.96
$ws.Cells.Item(132, 12).Value = 14705.3568
$ws.Cells.Item(132, 13).Value = -1328.96
$ws.Cells.Item(132, 14).Value = -19765.3568

$ws = $wb.Worksheets.Item("BSM")
# Row 108
$ws.Cells.Item(108, 8).Value = 0
$ws.Cells.Item(108, 10).Value = 0
$ws.Cells.Item(108, 12).Value = 0
$ws.Cells.Item(108, 14).ClearContents()

# Row 134
$ws.Cells.Item(134, 8).Value = 2443.611
$ws.Cells.Item(134, 9).Value = 1230.9166
$ws.Cells.Item(134, 11).Value = 3692.7498
$ws.Cells.Item(134, 13).Value = -1157.7498

$ws = $wb.Worksheets.Item("CRP")
# Row 6
$ws.Cells.Item(6, 8).Value = 7399214.5
$ws.Cells.Item(6, 9).Value = 10358001
$ws.Cells.Item(6, 10).Value = 2250
$ws.Cells.Item(6, 11).Value = 10358001
$ws.Cells.Item(6, 12).Value = 2250
$ws.Cells.Item(6, 13).Value = -10357888
$ws.Cells.Item(6, 14).Value = -2476

# Row 31
$ws.Cells.Item(31, 8).Value = 3950.1943
$ws.Cells.Item(31, 9).Value = 892.76666
$ws.Cells.Item(31, 10).Value = 6134.0713
$ws.Cells.Item(31, 11).Value = 892.76666
$ws.Cells.Item(31, 12).Value = 6134.0713
$ws.Cells.Item(31, 13).Value = -597.76666
$ws.Cells.Item(31, 14).Value = -6724.0713

# Row 34
$ws.Cells.Item(34, 8).Value = 3950.1943
$ws.Cells.Item(34, 9).Value = 892.76666
$ws.Cells.Item(34, 10).Value = 6134.0713
$ws.Cells.Item(34, 11).Value = 892.76666
$ws.Cells.Item(34, 12).Value = 6134.0713
$ws.Cells.Item(34, 13).Value = -690.76666
$ws.Cells.Item(34, 14).Value = -6538.0713

# Row 74
$ws.Cells.Item(74, 8).Value = 18562.857
$ws.Cells.Item(74, 10).Value = 18562.857
$ws.Cells.Item(74, 12).Value = 18562.857
$ws.Cells.Item(74, 14).Value = -20310.857

# Row 77
$ws.Cells.Item(77, 8).Value = 18562.857
$ws.Cells.Item(77, 10).Value = 18562.857
$ws.Cells.Item(77, 12).Value = 55688.571
$ws.Cells.Item(77, 14).Value = -64424.571

# Row 99
$ws.Cells.Item(99, 8).Value = 4588133.5
$ws.Cells.Item(99, 9).Value = 6402887
$ws.Cells.Item(99, 11).Value = 6402887
$ws.Cells.Item(99, 13).Value = -6401389

# Row 126
$ws.Cells.Item(126, 8).Value = 4588133.5
$ws.Cells.Item(126, 9).Value = 6402887
$ws.Cells.Item(126, 11).Value = 19208661
$ws.Cells.Item(126, 13).Value = -19206191

$ws = $wb.Worksheets.Item("CUL")
# Row 68
$ws.Cells.Item(68, 8).Value = 564019.5
$ws.Cells.Item(68, 9).Value = 1488012.6
$ws.Cells.Item(68, 10).Value = 1588.9348
$ws.Cells.Item(68, 11).Value = 4464037.800000001
$ws.Cells.Item(68, 12).Value = 4766.8044
$ws.Cells.Item(68, 13).Value = -4463226.800000001
$ws.Cells.Item(68, 14).Value = -6388.8044

# Row 71
$ws.Cells.Item(71, 8).Value = 564019.5
$ws.Cells.Item(71, 9).Value = 1488012.6
$ws.Cells.Item(71, 10).Value = 1588.9348
$ws.Cells.Item(71, 11).Value = 13392113.4
$ws.Cells.Item(71, 12).Value = 14300.4132
$ws.Cells.Item(71, 13).Value = -13388057.4
$ws.Cells.Item(71, 14).Value = -22412.4132

# Row 98
$ws.Cells.Item(98, 8).Value = 499.08334
$ws.Cells.Item(98, 10).Value = 488.5
$ws.Cells.Item(98, 12).Value = 1465.5
$ws.Cells.Item(98, 14).Value = -4461.5

$ws = $wb.Worksheets.Item("GSM")
# Row 70
$ws.Cells.Item(70, 8).Value = 16084.75
$ws.Cells.Item(70, 9).Value = 23337.6
$ws.Cells.Item(70, 11).Value = 23337.6
$ws.Cells.Item(70, 13).Value = -23067.6

# Row 73
$ws.Cells.Item(73, 8).Value = 16084.75
$ws.Cells.Item(73, 9).Value = 23337.6
$ws.Cells.Item(73, 11).Value = 23337.6
$ws.Cells.Item(73, 13).Value = -22401.6

# Row 102
$ws.Cells.Item(102, 8).Value = 750.5714
$ws.Cells.Item(102, 9).Value = 550.6667
$ws.Cells.Item(102, 10).Value = 1950
$ws.Cells.Item(102, 11).Value = 550.6667
$ws.Cells.Item(102, 12).Value = 1950
$ws.Cells.Item(102, 13).Value = 1071.3333
$ws.Cells.Item(102, 14).Value = -5194

# Row 108
$ws.Cells.Item(108, 8).Value = 0
$ws.Cells.Item(108, 10).Value = 0
$ws.Cells.Item(108, 12).Value = 0
$ws.Cells.Item(108, 14).ClearContents()

# Row 135
$ws.Cells.Item(135, 8).Value = 43692.5
$ws.Cells.Item(135, 10).Value = 43692.5
$ws.Cells.Item(135, 12).Value = 43692.5
$ws.Cells.Item(135, 14).Value = -53832.5

$ws = $wb.Worksheets.Item("LTW")
# Row 22
$ws.Cells.Item(22, 8).Value = 999.2381
$ws.Cells.Item(22, 9).Value = 598.6667
$ws.Cells.Item(22, 10).Value = 2000.6666
$ws.Cells.Item(22, 11).Value = 598.6667
$ws.Cells.Item(22, 12).Value = 2000.6666
$ws.Cells.Item(22, 13).Value = -303.6667
$ws.Cells.Item(22, 14).Value = -2590.6666

# Row 27
$ws.Cells.Item(27, 8).Value = 999.2381
$ws.Cells.Item(27, 9).Value = 598.6667
$ws.Cells.Item(27, 10).Value = 2000.6666
$ws.Cells.Item(27, 11).Value = 598.6667
$ws.Cells.Item(27, 12).Value = 2000.6666
$ws.Cells.Item(27, 13).Value = -491.6667
$ws.Cells.Item(27, 14).Value = -2214.6666

# Row 46
$ws.Cells.Item(46, 8).Value = 72675.5
$ws.Cells.Item(46, 9).Value = 112219.78
$ws.Cells.Item(46, 11).Value = 112219.78
$ws.Cells.Item(46, 13).Value = -112031.78

# Row 93
$ws.Cells.Item(93, 8).Value = 9780.546
$ws.Cells.Item(93, 9).Value = 11760.667
$ws.Cells.Item(93, 11).Value = 11760.667
$ws.Cells.Item(93, 13).Value = -10512.667

# Row 122
$ws.Cells.Item(122, 8).Value = 27779448
$ws.Cells.Item(122, 9).Value = 37038372
$ws.Cells.Item(122, 10).Value = 2680
$ws.Cells.Item(122, 11).Value = 111115116
$ws.Cells.Item(122, 12).Value = 8040
$ws.Cells.Item(122, 13).Value = -111112666
$ws.Cells.Item(122, 14).Value = -12940

# Row 132
$ws.Cells.Item(132, 8).Value = 4545.4224
$ws.Cells.Item(132, 9).Value = 5168.6523
$ws.Cells.Item(132, 11).Value = 15505.9569
$ws.Cells.Item(132, 13).Value = -12975.9569

$ws = $wb.Worksheets.Item("WVR")
# Row 81
$ws.Cells.Item(81, 8).Value = 853.4666999999999
$ws.Cells.Item(81, 9).Value = 640.2
$ws.Cells.Item(81, 11).Value = 1280.4
$ws.Cells.Item(81, 13).Value = -219.4000000000001

# Row 84
$ws.Cells.Item(84, 8).Value = 853.4666999999999
$ws.Cells.Item(84, 9).Value = 640.2
$ws.Cells.Item(84, 11).Value = 6402
$ws.Cells.Item(84, 13).Value = -1098

# Row 96
$ws.Cells.Item(96, 8).Value = 1390
$ws.Cells.Item(96, 9).Value = 1202.4
$ws.Cells.Item(96, 10).Value = 1702.6666
$ws.Cells.Item(96, 11).Value = 1202.4
$ws.Cells.Item(96, 12).Value = 1702.6666
$ws.Cells.Item(96, 13).Value = 170.5999999999999
$ws.Cells.Item(96, 14).Value = -4448.6666

# Row 108
$ws.Cells.Item(108, 8).Value = 0
$ws.Cells.Item(108, 10).Value = 0
$ws.Cells.Item(108, 12).Value = 0
$ws.Cells.Item(108, 14).ClearContents()
